$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182277083396912
$ws.Range("B1").Value = 2.404687643051147
$ws.Range("C1").Value = 3.772290468215942
$ws.Range("D1").Value = 2.091436624526978
$ws.Range("E1").Value = 1.201637148857117
